$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.803.19"
$ws.Range("E2").Value = "  +0.68%  "
$ws.Range("D3").Value = "2.981.28"
$ws.Range("E3").Value = "  -0.04%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "588.53"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.36%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.91"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.25%  "
$ws.Range("E7").Value = "  +0.18%  "
$ws.Range("D8").Value = "2.966.07"
$ws.Range("E8").Value = "  -0.49%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.503"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.78%  "
$ws.Range("E10").Value = "  +11.57%  "
$ws.Range("E11").Value = "  -1.24%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.454"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.51%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000226"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.39%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.66"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.77%  "
$ws.Range("E15").Value = "  -0.73%  "
$ws.Range("D16").Value = "3.476.53"
$ws.Range("E16").Value = "  +0.07%  "
$ws.Range("D17").Value = "61.907.15"
$ws.Range("E17").Value = "  +0.86%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.93"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.38%  "
$ws.Range("D19").Value = "2.980.62"
$ws.Range("E19").Value = "  -0.11%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "438.74"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.15%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.96"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.49%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.682"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.16%  "
$ws.Range("E23").Value = "  +0.99%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "81.51"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.75%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.10"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.21%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.21"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.69%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.99"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.28%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.33"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.23%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.22"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +8.35%  "
$ws.Range("B31").Value = "FirstDigitalUSD"
$ws.Range("C31").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.00"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.08%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.65"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.73%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.07"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.03%  "
$ws.Range("E34").Value = "  -1.20%  "
$ws.Range("D35").Value = "0.0₃0841"
$ws.Range("E35").Value = "  +3.95%  "
$ws.Range("E36").Value = "  +0.86%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.75"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.02%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.01"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.97%  "
$ws.Range("E39").Value = "  -0.22%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.04"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.61%  "
$ws.Range("E41").Value = "  +1.09%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.82"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.34%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.297"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +10.65%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "42.72"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +8.04%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0351"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.90%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "372.09"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.84%  "
$ws.Range("D47").Value = "2.673.44"
$ws.Range("E47").Value = "  -0.40%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "132.18"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.54%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "25.78"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +10.54%  "
$ws.Range("E50").Value = "  +0.00%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.20"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.40%  "
